$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original content of every data row (2-27) before overwriting anything,
# since the edit is a permutation (rows get re-ordered / re-assigned among each other).
$row2 = $ws.Range("A2:R2").Value2
$row3 = $ws.Range("A3:R3").Value2
$row4 = $ws.Range("A4:R4").Value2
$row5 = $ws.Range("A5:R5").Value2
$row6 = $ws.Range("A6:R6").Value2
$row7 = $ws.Range("A7:R7").Value2
$row8 = $ws.Range("A8:R8").Value2
$row9 = $ws.Range("A9:R9").Value2
$row10 = $ws.Range("A10:R10").Value2
$row11 = $ws.Range("A11:R11").Value2
$row12 = $ws.Range("A12:R12").Value2
$row13 = $ws.Range("A13:R13").Value2
$row14 = $ws.Range("A14:R14").Value2
$row15 = $ws.Range("A15:R15").Value2
$row16 = $ws.Range("A16:R16").Value2
$row17 = $ws.Range("A17:R17").Value2
$row18 = $ws.Range("A18:R18").Value2
$row19 = $ws.Range("A19:R19").Value2
$row20 = $ws.Range("A20:R20").Value2
$row21 = $ws.Range("A21:R21").Value2
$row22 = $ws.Range("A22:R22").Value2
$row23 = $ws.Range("A23:R23").Value2
$row24 = $ws.Range("A24:R24").Value2
$row25 = $ws.Range("A25:R25").Value2
$row26 = $ws.Range("A26:R26").Value2
$row27 = $ws.Range("A27:R27").Value2

# Write back rows in their new order per the target permutation.
$ws.Range("A2:R2").Value = $row16
$ws.Range("A3:R3").Value = $row17
$ws.Range("A4:R4").Value = $row15
$ws.Range("A5:R5").Value = $row24
$ws.Range("A6:R6").Value = $row5
$ws.Range("A7:R7").Value = $row18
$ws.Range("A8:R8").Value = $row26
$ws.Range("A9:R9").Value = $row27
$ws.Range("A10:R10").Value = $row13
$ws.Range("A11:R11").Value = $row19
$ws.Range("A12:R12").Value = $row11
$ws.Range("A13:R13").Value = $row4
$ws.Range("A14:R14").Value = $row10
$ws.Range("A15:R15").Value = $row12
$ws.Range("A16:R16").Value = $row7
$ws.Range("A17:R17").Value = $row20
$ws.Range("A18:R18").Value = $row25
$ws.Range("A19:R19").Value = $row14
$ws.Range("A20:R20").Value = $row9
$ws.Range("A21:R21").Value = $row23
$ws.Range("A22:R22").Value = $row21
$ws.Range("A23:R23").Value = $row22
$ws.Range("A24:R24").Value = $row2
$ws.Range("A25:R25").Value = $row3
$ws.Range("A26:R26").Value = $row6
$ws.Range("A27:R27").Value = $row8
